$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5-25 shift down to 6-26.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with data (mirrors the row above it, with the
# changes described in the diff).
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44953
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103002
$ws.Range("J5").Value = "Ciruela"
$ws.Range("K5").Value = "Black Amber"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 350
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19571
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1087
$ws.Range("T5").Value = 18
